$d = $word.ActiveDocument

$replacements = @(
    @("61÷2=30, 1", "87÷2=43, 1"),
    @("66÷3=22, 0", "30÷3=10, 0"),
    @("97÷7=13, 6", "42÷9=4, 6"),
    @("49÷4=12, 1", "91÷5=18, 1"),
    @("76÷7=10, 6", "29÷6=4, 5"),
    @("48÷2=24, 0", "16÷9=1, 7"),
    @("60÷5=12, 0", "93÷8=11, 5"),
    @("41÷9=4, 5", "14÷5=2, 4"),
    @("38÷6=6, 2", "99÷7=14, 1"),
    @("24÷8=3, 0", "98÷2=49, 0"),
    @("14÷2=7, 0", "22÷2=11, 0"),
    @("95÷6=15, 5", "51÷6=8, 3"),
    @("27÷9=3, 0", "42÷6=7, 0"),
    @("34÷6=5, 4", "89÷8=11, 1"),
    @("76÷3=25, 1", "33÷3=11, 0"),
    @("81÷3=27, 0", "81÷9=9, 0"),
    @("12÷6=2, 0", "44÷7=6, 2"),
    @("18÷3=6, 0", "32÷8=4, 0"),
    @("17÷3=5, 2", "84÷3=28, 0"),
    @("75÷2=37, 1", "83÷6=13, 5"),
    @("12÷7=1, 5", "43÷3=14, 1"),
    @("58÷2=29, 0", "69÷3=23, 0"),
    @("58÷6=9, 4", "70÷4=17, 2"),
    @("70÷6=11, 4", "81÷9=9, 0"),
    @("45÷6=7, 3", "67÷6=11, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
